{"js": "// Applies the German copy-edit pass described in the commit diff.\n// Each entry is an exact, unique source substring paired with its\n// replacement. We use Word's search API (exact match, case sensitive,\n// punctuation/whitespace respected) and replace each hit in place via\n// Range.insertText(..., Word.InsertLocation.replace) so the existing\n// run formatting (rPr) on the hit is preserved, matching how Word's\n// own Find & Replace keeps formatting when the match falls inside a\n// single run.\nconst replacements = [\n  [\n    \" Die virale Rolle vermittelte erfolgreich den Reiz des Produkts durch ansprechende visuelle und informative Inhalte.\",\n    \" Das virale Reel vermittelte erfolgreich die Attraktivit\u00e4t des Produkts durch ansprechende Visuals und informative Inhalte.\"\n  ],\n  [\n    \"Influencer Marketing:\",\n    \"Influencer-Marketing:\"\n  ],\n  [\n    \" Die Macht des Influencer-Marketings kann nicht \u00fcberstatiert werden.\",\n    \" Die Macht des Influencer-Marketings kann gar nicht hoch genug eingesch\u00e4tzt werden.\"\n  ],\n  [\n    \"Geschmacks- und Geschmackssorten:\",\n    \"Geschmack und Geschmackssorten:\"\n  ],\n  [\n    \" Der Ruf von Contoso Protein Plus f\u00fcr k\u00f6stliche und vielf\u00e4ltige Aromen war ein wichtiger Verkaufspunkt im Viral-Inhalt.\",\n    \" Der Ruf von Contoso Protein Plus, lecker und abwechslungsreich zu sein, war ein wichtiges Verkaufsargument in den viralen Inhalten.\"\n  ],\n  [\n    \" Der anhaltende Anstieg des Gesundheits- und Fitnessbewusstseins, kombiniert mit einer Zunahme der Anzahl der Menschen, die Trainingsroutinen und aktive Lebensstile annehmen, schuf einen empf\u00e4nglichen Markt f\u00fcr ein Produkt wie Contoso Protein Plus.\",\n    \" Das st\u00e4ndig wachsende Bewusstsein f\u00fcr Gesundheit und Fitness, verbunden mit einer steigenden Anzahl von Menschen, die eine Trainingsroutine und einen aktiven Lebensstil verfolgen, hat einen empf\u00e4nglichen Markt f\u00fcr ein Produkt wie Contoso Protein Plus geschaffen.\"\n  ],\n  [\n    \" Die Barrierefreiheit des Produkts durch verschiedene Onlineh\u00e4ndler hat den Hype weiter gef\u00f6rdert.\",\n    \" Die Verf\u00fcgbarkeit des Produkts \u00fcber verschiedene Onlineh\u00e4ndler hat den Hype weiter angeheizt.\"\n  ],\n  [\n    \" Die Rolle war kein Isolierter Fall.\",\n    \" Das Reel war kein Einzelfall.\"\n  ],\n  [\n    \"Wort des Mundes:\",\n    \"Mundpropaganda:\"\n  ],\n  [\n    \" Social Media-Plattformen f\u00f6rdern die rasante Verbreitung von Trends durch Mundwort.\",\n    \" Social Media-Plattformen f\u00f6rdern die schnelle Verbreitung von Trends durch Mundpropaganda.\"\n  ]\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, replacementText] of replacements) {\n  const results = body.search(searchText, {\n    matchCase: true,\n    matchWholeWord: false,\n    matchWildcards: false\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + searchText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replacementText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the German copy-edit pass described in the commit diff.\n# Each pair is an exact, unique source substring and its replacement.\n# We drive Word's Find/Replace (Range.Find.Execute) against the whole\n# document story so the existing run formatting on each hit is kept,\n# matching how Word's own Find & Replace behaves when a match sits\n# inside a single run.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \" Die virale Rolle vermittelte erfolgreich den Reiz des Produkts durch ansprechende visuelle und informative Inhalte.\"; Replace = \" Das virale Reel vermittelte erfolgreich die Attraktivit\u00e4t des Produkts durch ansprechende Visuals und informative Inhalte.\" },\n    @{ Find = \"Influencer Marketing:\"; Replace = \"Influencer-Marketing:\" },\n    @{ Find = \" Die Macht des Influencer-Marketings kann nicht \u00fcberstatiert werden.\"; Replace = \" Die Macht des Influencer-Marketings kann gar nicht hoch genug eingesch\u00e4tzt werden.\" },\n    @{ Find = \"Geschmacks- und Geschmackssorten:\"; Replace = \"Geschmack und Geschmackssorten:\" },\n    @{ Find = \" Der Ruf von Contoso Protein Plus f\u00fcr k\u00f6stliche und vielf\u00e4ltige Aromen war ein wichtiger Verkaufspunkt im Viral-Inhalt.\"; Replace = \" Der Ruf von Contoso Protein Plus, lecker und abwechslungsreich zu sein, war ein wichtiges Verkaufsargument in den viralen Inhalten.\" },\n    @{ Find = \" Der anhaltende Anstieg des Gesundheits- und Fitnessbewusstseins, kombiniert mit einer Zunahme der Anzahl der Menschen, die Trainingsroutinen und aktive Lebensstile annehmen, schuf einen empf\u00e4nglichen Markt f\u00fcr ein Produkt wie Contoso Protein Plus.\"; Replace = \" Das st\u00e4ndig wachsende Bewusstsein f\u00fcr Gesundheit und Fitness, verbunden mit einer steigenden Anzahl von Menschen, die eine Trainingsroutine und einen aktiven Lebensstil verfolgen, hat einen empf\u00e4nglichen Markt f\u00fcr ein Produkt wie Contoso Protein Plus geschaffen.\" },\n    @{ Find = \" Die Barrierefreiheit des Produkts durch verschiedene Onlineh\u00e4ndler hat den Hype weiter gef\u00f6rdert.\"; Replace = \" Die Verf\u00fcgbarkeit des Produkts \u00fcber verschiedene Onlineh\u00e4ndler hat den Hype weiter angeheizt.\" },\n    @{ Find = \" Die Rolle war kein Isolierter Fall.\"; Replace = \" Das Reel war kein Einzelfall.\" },\n    @{ Find = \"Wort des Mundes:\"; Replace = \"Mundpropaganda:\" },\n    @{ Find = \" Social Media-Plattformen f\u00f6rdern die rasante Verbreitung von Trends durch Mundwort.\"; Replace = \" Social Media-Plattformen f\u00f6rdern die schnelle Verbreitung von Trends durch Mundpropaganda.\" }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n    if (-not $found) {\n        throw \"No match found for: $($pair.Find)\"\n    }\n}\n"}
